$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 574, shifting existing rows 574:617 down to 575:618
$ws.Rows.Item(574).Insert()

# Populate the new row 574 with the weekly data point (same Variedad/Calidad/Origen as old row 574, new date/prices)
$ws.Range("A574").Value = 8
$ws.Range("B574").Value = "Terminal La Palmera de La Serena"
$ws.Range("C574").Value = "Coquimbo"
$ws.Range("D574").Value = 45021
$ws.Range("E574").Value = 4
$ws.Range("F574").Value = 100114001
$ws.Range("G574").Value = "Papa"
$ws.Range("H574").Value = "Cardinal"
$ws.Range("I574").Value = "1a (cosecha)"
$ws.Range("J574").Value = 2000
$ws.Range("K574").Value = 12500
$ws.Range("L574").Value = 13000
$ws.Range("M574").Value = 12750
$ws.Range("N574").Value = "`$/saco 25 kilos"
$ws.Range("O574").Value = "Provincia del Elquí"
$ws.Range("P574").Value = 510
$ws.Range("Q574").Value = 25
$ws.Range("R574").Value = "Hortaliza"
